$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("novosjogos (2)")

# Row 5: Situacao -> Concluido, Data Termino -> 09/01/2025 (date instead of text)
$ws.Range("D5").Value = "Concluido"
$ws.Range("E5").Value = [DateTime]"2025-01-09"
$ws.Range("E5").NumberFormat = "dd/mm/yyyy"

# Row 6: Situacao -> Concluido (keep existing Data Termino)
$ws.Range("D6").Value = "Concluido"

# Update selection to D6
$ws.Range("D6").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
